# Add "Checked" columns (INV.no Check, Inv. Date Check, ExcludeVAT_diff,
# VAT_diff, IncludeVAT_diff) to the reconciliation sheets, pushing the
# existing "null_report" column out to the right, and update the summary
# "report" sheet to drop the now-obsolete "DIFF" bucket.

$wb = $excel.ActiveWorkbook

function Add-CheckColumns($ws, $lastRow) {
    # Column layout (1-based): F/G = BASE invoice no/date, H/I/J = BASE
    # amounts, Q/R = CPFM invoice no/date, S/T/U = CPFM amounts.
    # V = null_report (existing) -> becomes AA; new V..Z inserted before it.

    # Copy the header formatting from the old "null_report" header (V1) onto
    # the five new header cells so they pick up the same style (bold,
    # bordered, centered) instead of default formatting.
    $ws.Range("V1").Copy() | Out-Null
    $ws.Range("W1:AA1").PasteSpecial(-4122) | Out-Null

    $ws.Range("V1").Value2 = "INV.no Check"
    $ws.Range("W1").Value2 = "Inv. Date Check"
    $ws.Range("X1").Value2 = "ExcludeVAT_diff"
    $ws.Range("Y1").Value2 = "VAT_diff"
    $ws.Range("Z1").Value2 = "IncludeVAT_diff"
    $ws.Range("AA1").Value2 = "null_report"

    for ($r = 2; $r -le $lastRow; $r++) {
        # Stash the old null_report value (currently sitting in column V)
        # before it gets overwritten by the new INV.no Check column.
        $nullReport = $ws.Cells.Item($r, 22).Value2

        $invBase = $ws.Cells.Item($r, 6).Value2
        $dateBase = $ws.Cells.Item($r, 7).Value2
        $excVat = $ws.Cells.Item($r, 8).Value2
        $taxBase = $ws.Cells.Item($r, 9).Value2
        $totalBase = $ws.Cells.Item($r, 10).Value2

        $invCpfm = $ws.Cells.Item($r, 17).Value2
        $dateCpfm = $ws.Cells.Item($r, 18).Value2
        $sumNett = $ws.Cells.Item($r, 19).Value2
        $taxCpfm = $ws.Cells.Item($r, 20).Value2
        $totalCpfm = $ws.Cells.Item($r, 21).Value2

        $baseMissing = ($invBase -eq $null -or $invBase -eq "")
        $cpfmMissing = ($invCpfm -eq $null -or $invCpfm -eq "")

        if ($baseMissing -or $cpfmMissing) {
            $ws.Cells.Item($r, 22).Value2 = $false
            $ws.Cells.Item($r, 23).Value2 = $false
            $ws.Cells.Item($r, 24).Value2 = ""
            $ws.Cells.Item($r, 25).Value2 = ""
            $ws.Cells.Item($r, 26).Value2 = ""
        }
        else {
            $ws.Cells.Item($r, 22).Value2 = [bool]($invBase -eq $invCpfm)
            $ws.Cells.Item($r, 23).Value2 = [bool]($dateBase -eq $dateCpfm)
            $ws.Cells.Item($r, 24).Value2 = [Math]::Round($excVat - $sumNett, 2)
            $ws.Cells.Item($r, 25).Value2 = [Math]::Round($taxBase - $taxCpfm, 2)
            $ws.Cells.Item($r, 26).Value2 = [Math]::Round($totalBase - $totalCpfm, 2)
        }

        if ($nullReport -eq $null -or $nullReport -eq "") {
            $ws.Cells.Item($r, 27).Value2 = ""
        }
        else {
            $ws.Cells.Item($r, 27).Value2 = $nullReport
        }
    }
}

# --- Reconciled Data (84 rows incl. header) ---
$wsMain = $wb.Worksheets.Item("Reconciled Data")
Add-CheckColumns $wsMain 84

# --- BASE_null (11 rows incl. header) ---
$wsBaseNull = $wb.Worksheets.Item("BASE_null")
Add-CheckColumns $wsBaseNull 11

# --- CPFM_null (6 rows incl. header) ---
$wsCpfmNull = $wb.Worksheets.Item("CPFM_null")
Add-CheckColumns $wsCpfmNull 6

# --- report: drop the obsolete "DIFF" row (the 5 rows it counted now live
# inside Reconciled Data with a populated check/diff instead of a separate
# bucket), shifting "Matching"/"Total" up and updating the Total.
$wsReport = $wb.Worksheets.Item("report")
$wsReport.Rows.Item(3).Delete() | Out-Null
$wsReport.Range("B4").Value2 = 83
